$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data changes (matches shared-string insertion order of the target) ---

# K2: "Test1" -> "test"
$ws.Range("K2").Value = 'test'

# A2: "Shanawaz" -> quote-prefixed text "123"
$ws.Range("A2").Value = "'123"

# S2:V2 misc text
$ws.Range("S2").Value = '##'
$ws.Range("T2").Value = '$$'
$ws.Range("U2").Value = '!!'
$ws.Range("V2").Value = 'Test11$'

# N2:R2 quote-prefixed decimal-looking text
$ws.Range("N2").Value = "'1.1"
$ws.Range("O2").Value = "'2.1"
$ws.Range("P2").Value = "'3.1"
$ws.Range("Q2").Value = "'4.4"
$ws.Range("R2").Value = "'5.1"

# --- Row 1 header text changes ---
$ws.Range("B1").Value = 'From Date'
$ws.Range("C1").Value = 'To Date'

# B2 / C2: new dates, with the updated custom date/time format
$ws.Range("B2").Value = "'02-Apr-2018  10:20"
$ws.Range("C2").Value = "'05-Apr-2018  10:21"
$ws.Range("B2").NumberFormat = 'dd\-mmm\-yyyy\ hh:mm'
$ws.Range("C2").NumberFormat = 'dd\-mmm\-yyyy\ hh:mm'

# --- Column widths / layout ---
$ws.Columns.Item(1).ColumnWidth = 12.333333333333334
$ws.Columns.Item(2).ColumnWidth = 16.333333333333332
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668

# --- Selection moves from B2 to C2 ---
$ws.Range("C2").Select() | Out-Null
